$wb = $excel.ActiveWorkbook

# The sheet that was named "Sheet1" is renamed to the Persian month name
# "فروردین" (Farvardin). It is the workbook's only/active sheet.
$ws = $wb.ActiveSheet
$ws.Name = "فروردین"
